# Auto-generated: updates currentAveragePrice / LevePrice / LeveProfit
# columns (H-N) for specific leve rows across multiple sheets,
# matching the scheduled market-data refresh described in the commit.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 4: Root Rush
$ws.Range("H4").Value = 4177.1
$ws.Range("I4").Value = 254.2
$ws.Range("J4").Value = 8100
$ws.Range("K4").Value = 254.2
$ws.Range("L4").Value = 8100
$ws.Range("M4").Value = -140.2
$ws.Range("N4").Value = -8328
# Row 17: One for the Road
$ws.Range("H17").Value = 3027023.8
$ws.Range("J17").Value = 3027023.8
$ws.Range("L17").Value = 9081071.399999999
$ws.Range("N17").Value = -9081407.399999999
# Row 40: Stuck in the Moment
$ws.Range("H40").Value = 1523.3646
$ws.Range("I40").Value = 1495.8148
$ws.Range("J40").Value = 1672.1333
$ws.Range("K40").Value = 1495.8148
$ws.Range("L40").Value = 1672.1333
$ws.Range("M40").Value = -1320.8148
$ws.Range("N40").Value = -2022.1333
# Row 76: Warding Off Temptation
$ws.Range("H76").Value = 4394.067
$ws.Range("I76").Value = 3271.8572
$ws.Range("J76").Value = 5376
$ws.Range("K76").Value = 3271.8572
$ws.Range("L76").Value = 5376
$ws.Range("M76").Value = -2956.8572
$ws.Range("N76").Value = -6006
# Row 79: The Garden of Arcane Delights (L)
$ws.Range("H79").Value = 4394.067
$ws.Range("I79").Value = 3271.8572
$ws.Range("J79").Value = 5376
$ws.Range("K79").Value = 3271.8572
$ws.Range("L79").Value = 5376
$ws.Range("M79").Value = -2179.8572
$ws.Range("N79").Value = -7560
# Row 100: Asking for a Friend
$ws.Range("H100").Value = 8334962.5
$ws.Range("I100").Value = 23810966
$ws.Range("J100").Value = 1729.3846
$ws.Range("K100").Value = 23810966
$ws.Range("L100").Value = 1729.3846
$ws.Range("M100").Value = -23810425
$ws.Range("N100").Value = -2811.3846
# Row 138: All-night Crafting
$ws.Range("H138").Value = 4000
$ws.Range("J138").Value = 8554.25
$ws.Range("L138").Value = 25662.75
$ws.Range("N138").Value = -35942.75
# Row 141: Remedy for Reason
$ws.Range("H141").Value = 2555.2273
$ws.Range("I141").Value = 2006.1111
$ws.Range("J141").Value = 5026.25
$ws.Range("K141").Value = 6018.3333
$ws.Range("L141").Value = 15078.75
$ws.Range("M141").Value = -838.3333000000002
$ws.Range("N141").Value = -25438.75

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 6898.2383
$ws.Range("I61").Value = 8258.625
$ws.Range("J61").Value = 2545
$ws.Range("K61").Value = 8258.625
$ws.Range("L61").Value = 2545
$ws.Range("M61").Value = -8046.625
$ws.Range("N61").Value = -2969
# Row 74: As the Bolt Flies
$ws.Range("H74").Value = 1457.1562
$ws.Range("I74").Value = 1379.965
$ws.Range("K74").Value = 1379.965
$ws.Range("M74").Value = -505.9649999999999
# Row 77: Heavy Metal Banned (L)
$ws.Range("H77").Value = 1457.1562
$ws.Range("I77").Value = 1379.965
$ws.Range("K77").Value = 6899.825
$ws.Range("M77").Value = -2531.825
# Row 102: Smells of Rich Tama-hagane
$ws.Range("H102").Value = 7410645
$ws.Range("I102").Value = 9262056
$ws.Range("J102").Value = 5000
$ws.Range("K102").Value = 9262056
$ws.Range("L102").Value = 5000
$ws.Range("M102").Value = -9260434
$ws.Range("N102").Value = -8244
# Row 122: Haste for High Durium
$ws.Range("H122").Value = 1711343.4
$ws.Range("I122").Value = 2138354.2
$ws.Range("J122").Value = 3300
$ws.Range("K122").Value = 6415062.600000001
$ws.Range("L122").Value = 9900
$ws.Range("M122").Value = -6412612.600000001
$ws.Range("N122").Value = -14800
# Row 123: The Armoire Is Open
$ws.Range("H123").Value = 40426
$ws.Range("J123").Value = 40426
$ws.Range("L123").Value = 40426
$ws.Range("N123").Value = -50226
# Row 136: Metal with Mettle
$ws.Range("H136").Value = 6898.2383
$ws.Range("I136").Value = 8258.625
$ws.Range("J136").Value = 2545
$ws.Range("K136").Value = 24775.875
$ws.Range("L136").Value = 7635
$ws.Range("M136").Value = -22225.875
$ws.Range("N136").Value = -12735

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 99: Meddle in Metal
$ws.Range("H99").Value = 76924350
$ws.Range("I99").Value = 100001000
$ws.Range("J99").Value = 2200
$ws.Range("K99").Value = 100001000
$ws.Range("L99").Value = 2200
$ws.Range("M99").Value = -99999502
$ws.Range("N99").Value = -5196
# Row 105: Ingot to Wing It
$ws.Range("H105").Value = 7914.6113
$ws.Range("I105").Value = 11923.9
$ws.Range("J105").Value = 2903
$ws.Range("K105").Value = 11923.9
$ws.Range("L105").Value = 2903
$ws.Range("M105").Value = -10176.9
$ws.Range("N105").Value = -6397
# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 4527.2197
$ws.Range("I134").Value = 5114.067
$ws.Range("J134").Value = 2926.7273
$ws.Range("K134").Value = 15342.201
$ws.Range("L134").Value = 8780.1819
$ws.Range("M134").Value = -12807.201
$ws.Range("N134").Value = -13850.1819

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 105: Zelkova, My Love
$ws.Range("H105").Value = 2287
$ws.Range("I105").Value = 2822
$ws.Range("J105").Value = 949.5
$ws.Range("K105").Value = 2822
$ws.Range("L105").Value = 949.5
$ws.Range("M105").Value = -1075
$ws.Range("N105").Value = -4443.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 122: Awarding Academic Excellence
$ws.Range("H122").Value = 2092372.9
$ws.Range("I122").Value = 2702042
$ws.Range("J122").Value = 2078.8572
$ws.Range("K122").Value = 8106126
$ws.Range("L122").Value = 6236.571599999999
$ws.Range("M122").Value = -8103676
$ws.Range("N122").Value = -11136.5716

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 46: Supply Side Logic
$ws.Range("H46").Value = 37037892
$ws.Range("I46").Value = 66667220
$ws.Range("J46").Value = 1237.5
$ws.Range("K46").Value = 66667220
$ws.Range("L46").Value = 1237.5
$ws.Range("M46").Value = -66667032
$ws.Range("N46").Value = -1613.5
# Row 122: Hell on Leather
$ws.Range("H122").Value = 5497106
$ws.Range("J122").Value = 2450
$ws.Range("L122").Value = 7350
$ws.Range("N122").Value = -12250

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 33: I'll Be Your Wailer Today
$ws.Range("H33").Value = 23250
$ws.Range("J33").Value = 23250
$ws.Range("L33").Value = 23250
$ws.Range("N33").Value = -23750
# Row 36: Put a Lid on It
$ws.Range("H36").Value = 23250
$ws.Range("J36").Value = 23250
$ws.Range("L36").Value = 23250
$ws.Range("N36").Value = -23750
# Row 37: Bet You Anything
$ws.Range("H37").Value = 23250
$ws.Range("J37").Value = 23250
$ws.Range("L37").Value = 23250
$ws.Range("N37").Value = -23656
# Row 75: Storm upon Bald Mountain
$ws.Range("H75").Value = 39150
$ws.Range("J75").Value = 39150
$ws.Range("L75").Value = 39150
$ws.Range("N75").Value = -41022
# Row 78: Abrupt Apprentices (L)
$ws.Range("H78").Value = 39150
$ws.Range("J78").Value = 39150
$ws.Range("L78").Value = 117450
$ws.Range("N78").Value = -126810
# Row 113: A Tender Table
$ws.Range("H113").Value = 1384.0834
$ws.Range("I113").Value = 2333
$ws.Range("J113").Value = 1067.7778
$ws.Range("K113").Value = 6999
$ws.Range("L113").Value = 3203.3334
$ws.Range("M113").Value = -4829
$ws.Range("N113").Value = -7543.3334
# Row 122: Heavy Armoire
$ws.Range("H122").Value = 1763.6364
$ws.Range("I122").Value = 2340
$ws.Range("J122").Value = 1283.3334
$ws.Range("K122").Value = 7020
$ws.Range("L122").Value = 3850.0002
$ws.Range("M122").Value = -4570
$ws.Range("N122").Value = -8750.0002
